$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.314.00"
$ws.Range("E2").Value = "  +1.55%  "
$ws.Range("D3").Value = "3.026.49"
$ws.Range("E3").Value = "  +1.19%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "519.33"
$ws.Range("E5").Value = "  +5.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.74"
$ws.Range("E6").Value = "  +5.66%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.439"
$ws.Range("E8").Value = "  +3.30%  "
$ws.Range("E9").Value = "  +5.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.109"
$ws.Range("E10").Value = "  +5.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.360"
$ws.Range("E11").Value = "  +2.15%  "
$ws.Range("E12").Value = "  +2.38%  "
$ws.Range("D13").Value = "3.543.97"
$ws.Range("E13").Value = "  +1.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.22"
$ws.Range("E14").Value = "  +5.16%  "
$ws.Range("E15").Value = "  +11.78%  "
$ws.Range("D16").Value = "57.308.04"
$ws.Range("E16").Value = "  +1.66%  "
$ws.Range("D17").Value = "3.022.55"
$ws.Range("E17").Value = "  +0.93%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.03"
$ws.Range("E18").Value = "  +2.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.73"
$ws.Range("E19").Value = "  +3.40%  "
$ws.Range("E20").Value = "  +3.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "332.74"
$ws.Range("E21").Value = "  +3.04%  "
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("E23").Value = "  +5.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.10"
$ws.Range("E24").Value = "  +5.06%  "
$ws.Range("E25").Value = "  +6.01%  "
$ws.Range("E26").Value = "  +0.31%  "
$ws.Range("D27").Value = "0.0₃0927"
$ws.Range("E27").Value = "  +5.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.80"
$ws.Range("E28").Value = "  +3.78%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.22"
$ws.Range("E29").Value = "  +7.23%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.83"
$ws.Range("E30").Value = "  +6.58%  "
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.23"
$ws.Range("E31").Value = "  +4.58%  "
$ws.Range("E32").Value = "  +4.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "158.67"
$ws.Range("E33").Value = "  +5.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.67"
$ws.Range("E34").Value = "  +4.23%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.80"
$ws.Range("E35").Value = "  +3.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.31"
$ws.Range("E36").Value = "  +2.43%  "
$ws.Range("B37").Value = "EnergySwap"
$ws.Range("C37").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "24.52"
$ws.Range("E37").Value = "  +3.57%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0679"
$ws.Range("E38").Value = "  +2.43%  "
$ws.Range("D39").Value = "3.057.88"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.47"
$ws.Range("E40").Value = "  +0.73%  "
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.655"
$ws.Range("E42").Value = "  +2.57%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.303.64"
$ws.Range("E43").Value = "  +6.18%  "
$ws.Range("E44").Value = "  +5.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.45"
$ws.Range("E45").Value = "  +2.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.02"
$ws.Range("E46").Value = "  +0.43%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.03"
$ws.Range("E47").Value = "  +10.08%  "
$ws.Range("E48").Value = "  +2.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.93"
$ws.Range("E49").Value = "  +6.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.60"
$ws.Range("E50").Value = "  +1.59%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0885"
$ws.Range("E51").Value = "  +4.18%  "
